$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells that receive numeric-looking text to stay as text,
# matching the original inlineStr cell type, then restore the default (unstyled) look.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.688.37"
$ws.Range("D3").Value = "2.356.92"
$ws.Range("D5").Value = "325.36"
$ws.Range("D6").Value = "100.88"
$ws.Range("D9").Value = "0.624"
$ws.Range("D10").Value = "40.00"
$ws.Range("D11").Value = "0.0923"
$ws.Range("D12").Value = "8.45"
$ws.Range("D15").Value = "16.51"
$ws.Range("D16").Value = "2.712.64"
$ws.Range("D17").Value = "2.353.00"
$ws.Range("D18").Value = "8.04"
$ws.Range("D19").Value = "42.650.79"
$ws.Range("D21").Value = "76.19"
$ws.Range("D23").Value = "266.20"
$ws.Range("D25").Value = "10.07"
$ws.Range("D27").Value = "11.46"
$ws.Range("D28").Value = "22.93"
$ws.Range("D30").Value = "175.89"
$ws.Range("D31").Value = "3.09"
$ws.Range("D32").Value = "0.0899"
$ws.Range("D33").Value = "35.26"
$ws.Range("D34").Value = "6.04"
$ws.Range("D37").Value = "0.0358"
$ws.Range("D38").Value = "2.93"
$ws.Range("D40").Value = "3.80"
$ws.Range("D41").Value = "1.51"
$ws.Range("D42").Value = "0.235"
$ws.Range("D43").Value = "69.90"
$ws.Range("D45").Value = "119.50"
$ws.Range("D46").Value = "90.75"
$ws.Range("D47").Value = "11.89"
$ws.Range("D48").Value = "5.52"
$ws.Range("D49").Value = "9.23"
$ws.Range("D50").Value = "1.26"
$ws.Range("D51").Value = "72.76"

$priceRange.Style = "Normal"

# Remaining text cells (Coin name, Link, Volume) — plain text, no coercion risk.
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("E6").Value = "  -8.38%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("E18").Value = "  +10.92%  "
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  +7.44%  "
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  -10.27%  "
$ws.Range("E25").Value = "  +9.24%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -10.24%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("E36").Value = "  -8.79%  "
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("E38").Value = "  +7.98%  "
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  -8.99%  "
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  +7.56%  "
$ws.Range("E47").Value = "  -7.26%  "
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("E51").Value = "  +3.09%  "
